$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 28573302
$ws.Range("I137").Value = 1199.3914
$ws.Range("J137").Value = 83336500
$ws.Range("K137").Value = 3598.1742
$ws.Range("L137").Value = 250009500
$ws.Range("M137").Value = -1048.1742
$ws.Range("N137").Value = -250014600
$ws.Range("H138").Value = 2491.236
$ws.Range("I138").Value = 2148.2
$ws.Range("J138").Value = 2920.0312
$ws.Range("K138").Value = 6444.599999999999
$ws.Range("L138").Value = 8760.0936
$ws.Range("M138").Value = -1304.599999999999
$ws.Range("N138").Value = -19040.0936
$ws.Range("H141").Value = 3128
$ws.Range("I141").Value = 1435
$ws.Range("J141").Value = 9900
$ws.Range("K141").Value = 4305
$ws.Range("L141").Value = 29700
$ws.Range("M141").Value = 875
$ws.Range("N141").Value = -40060

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6087.9
$ws.Range("I32").Value = 6048.384
$ws.Range("K32").Value = 6048.384
$ws.Range("M32").Value = -5761.384
$ws.Range("H61").Value = 4669.913
$ws.Range("I61").Value = 5496.923
$ws.Range("J61").Value = 3594.8
$ws.Range("K61").Value = 5496.923
$ws.Range("L61").Value = 3594.8
$ws.Range("M61").Value = -5284.923
$ws.Range("N61").Value = -4018.8
$ws.Range("H74").Value = 6337.129
$ws.Range("I74").Value = 1072.2106
$ws.Range("K74").Value = 1072.2106
$ws.Range("M74").Value = -198.2106000000001
$ws.Range("H77").Value = 6337.129
$ws.Range("I77").Value = 1072.2106
$ws.Range("K77").Value = 5361.053000000001
$ws.Range("M77").Value = -993.0530000000008
$ws.Range("H136").Value = 4669.913
$ws.Range("I136").Value = 5496.923
$ws.Range("J136").Value = 3594.8
$ws.Range("K136").Value = 16490.769
$ws.Range("L136").Value = 10784.4
$ws.Range("M136").Value = -13940.769
$ws.Range("N136").Value = -15884.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 187.36667
$ws.Range("I80").Value = 40.714287
$ws.Range("J80").Value = 232
$ws.Range("K80").Value = 40.714287
$ws.Range("L80").Value = 232
$ws.Range("M80").Value = 957.285713
$ws.Range("N80").Value = -2228
$ws.Range("H81").Value = 25566.666
$ws.Range("J81").Value = 25566.666
$ws.Range("L81").Value = 25566.666
$ws.Range("N81").Value = -27688.666
$ws.Range("H83").Value = 187.36667
$ws.Range("I83").Value = 40.714287
$ws.Range("J83").Value = 232
$ws.Range("K83").Value = 203.571435
$ws.Range("L83").Value = 1160
$ws.Range("M83").Value = 4788.428565
$ws.Range("N83").Value = -11144
$ws.Range("H84").Value = 25566.666
$ws.Range("J84").Value = 25566.666
$ws.Range("L84").Value = 76699.99800000001
$ws.Range("N84").Value = -87307.99800000001
$ws.Range("H134").Value = 93766.75
$ws.Range("J134").Value = 2749.1667
$ws.Range("L134").Value = 8247.500100000001
$ws.Range("N134").Value = -13317.5001
$ws.Range("H135").Value = 61826.668
$ws.Range("J135").Value = 61826.668
$ws.Range("L135").Value = 61826.668
$ws.Range("N135").Value = -71966.66800000001
$ws.Range("H138").Value = 69975
$ws.Range("J138").Value = 69975
$ws.Range("L138").Value = 69975
$ws.Range("N138").Value = -80255

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2337.0732
$ws.Range("I31").Value = 1034.9615
$ws.Range("J31").Value = 4594.067
$ws.Range("K31").Value = 1034.9615
$ws.Range("L31").Value = 4594.067
$ws.Range("M31").Value = -739.9614999999999
$ws.Range("N31").Value = -5184.067
$ws.Range("H34").Value = 2337.0732
$ws.Range("I34").Value = 1034.9615
$ws.Range("J34").Value = 4594.067
$ws.Range("K34").Value = 1034.9615
$ws.Range("L34").Value = 4594.067
$ws.Range("M34").Value = -832.9614999999999
$ws.Range("N34").Value = -4998.067
$ws.Range("H132").Value = 1652.98
$ws.Range("I132").Value = 1507.186
$ws.Range("J132").Value = 2548.5715
$ws.Range("K132").Value = 4521.558
$ws.Range("L132").Value = 7645.7145
$ws.Range("M132").Value = -1991.558
$ws.Range("N132").Value = -12705.7145

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 1000
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H103").Value = 3778355
$ws.Range("I103").Value = 4250530.5
$ws.Range("J103").Value = 950
$ws.Range("K103").Value = 12751591.5
$ws.Range("L103").Value = 2850
$ws.Range("M103").Value = -12750712.5
$ws.Range("N103").Value = -4608
$ws.Range("H110").Value = 3524.75
$ws.Range("J110").Value = 4250
$ws.Range("L110").Value = 12750
$ws.Range("N110").Value = -20930
$ws.Range("H114").Value = 1259.4584
$ws.Range("I114").Value = 515.5
$ws.Range("J114").Value = 1790.8572
$ws.Range("K114").Value = 1546.5
$ws.Range("L114").Value = 5372.571599999999
$ws.Range("M114").Value = 1707.5
$ws.Range("N114").Value = -11880.5716
$ws.Range("H131").Value = 2709.8948
$ws.Range("J131").Value = 1742.2963
$ws.Range("L131").Value = 5226.8889
$ws.Range("N131").Value = -15306.8889

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2050.6667
$ws.Range("I7").Value = 1880.8
$ws.Range("J7").Value = 2900
$ws.Range("K7").Value = 1880.8
$ws.Range("L7").Value = 2900
$ws.Range("M7").Value = -1768.8
$ws.Range("N7").Value = -3124
$ws.Range("H122").Value = 1875.2858
$ws.Range("I122").Value = 1683.4
$ws.Range("J122").Value = 2355
$ws.Range("K122").Value = 5050.200000000001
$ws.Range("L122").Value = 7065
$ws.Range("M122").Value = -2600.200000000001
$ws.Range("N122").Value = -11965
$ws.Range("H126").Value = 2050.6667
$ws.Range("I126").Value = 1880.8
$ws.Range("J126").Value = 2900
$ws.Range("K126").Value = 5642.4
$ws.Range("L126").Value = 8700
$ws.Range("M126").Value = -3172.4
$ws.Range("N126").Value = -13640
$ws.Range("H136").Value = 1468.8572
$ws.Range("I136").Value = 1124
$ws.Range("J136").Value = 2733.3333
$ws.Range("K136").Value = 3372
$ws.Range("L136").Value = 8199.999899999999
$ws.Range("M136").Value = -822
$ws.Range("N136").Value = -13299.9999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 27761752
$ws.Range("I2").Value = 10000000
$ws.Range("J2").Value = 33682336
$ws.Range("K2").Value = 10000000
$ws.Range("L2").Value = 33682336
$ws.Range("M2").Value = -9999888
$ws.Range("N2").Value = -33682560
$ws.Range("H107").Value = 653.1905
$ws.Range("I107").Value = 395.8
$ws.Range("J107").Value = 1296.6666
$ws.Range("K107").Value = 1187.4
$ws.Range("L107").Value = 3889.9998
$ws.Range("M107").Value = 732.5999999999999
$ws.Range("N107").Value = -7729.9998
$ws.Range("H122").Value = 3199.7144
$ws.Range("I122").Value = 2076.7778
$ws.Range("K122").Value = 6230.3334
$ws.Range("M122").Value = -3780.3334
$ws.Range("H136").Value = 11508.8
$ws.Range("I136").Value = 12168.637
$ws.Range("J136").Value = 6670
$ws.Range("K136").Value = 36505.911
$ws.Range("L136").Value = 20010
$ws.Range("M136").Value = -33955.911
$ws.Range("N136").Value = -25110
